$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H44").Value = 22262.5
$ws.Range("J44").Value = 22262.5
$ws.Range("L44").Value = 22262.5
$ws.Range("N44").Value = -23186.5

$ws.Range("H95").Value = 28000
$ws.Range("J95").Value = 28000
$ws.Range("L95").Value = 28000
$ws.Range("N95").Value = -33492

$ws.Range("H132").Value = 8774201
$ws.Range("I132").Value = 9011174
$ws.Range("K132").Value = 27033522
$ws.Range("M132").Value = -27030992

$ws.Range("H141").Value = 886
$ws.Range("I141").Value = 886
$ws.Range("K141").Value = 2658
$ws.Range("M141").Value = 2522

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 1000
$ws.Range("I3").Value = 1000
$ws.Range("K3").Value = 1000
$ws.Range("M3").Value = -885

$ws.Range("H8").Value = 5001502.5
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()

$ws.Range("H11").Value = 5002999.5
$ws.Range("I11").Value = 5002999.5
$ws.Range("K11").Value = 5002999.5
$ws.Range("M11").Value = -5002855.5

$ws.Range("H36").Value = 5026
$ws.Range("I36").Value = 5026
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 5026
$ws.Range("L36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("N36").Value = -4680

$ws.Range("H46").Value = 5835.75
$ws.Range("J46").Value = 4450.6665
$ws.Range("L46").Value = 4450.6665
$ws.Range("N46").Value = -5088.6665

$ws.Range("H94").Value = 22222
$ws.Range("J94").Value = 22222
$ws.Range("L94").Value = 22222
$ws.Range("N94").Value = -24024

$ws.Range("H95").Value = 12671.667
$ws.Range("J95").Value = 12671.667
$ws.Range("L95").Value = 12671.667
$ws.Range("N95").Value = -18163.667

$ws.Range("H96").Value = 18499.75
$ws.Range("J96").Value = 18499.75
$ws.Range("L96").Value = 18499.75
$ws.Range("N96").Value = -23991.75

$ws.Range("H132").Value = 3002.318
$ws.Range("I132").Value = 3265.5557
$ws.Range("J132").Value = 2820.077
$ws.Range("K132").Value = 9796.667099999999
$ws.Range("L132").Value = 8460.231
$ws.Range("M132").Value = -7266.667099999999
$ws.Range("N132").Value = -13520.231

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 20999
$ws.Range("J92").Value = 20999
$ws.Range("L92").Value = 20999
$ws.Range("N92").Value = -25991

$ws.Range("H94").Value = 62501490
$ws.Range("I94").Value = 62501490
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 62501490
$ws.Range("L94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -62501039

$ws.Range("H105").Value = 71430456
$ws.Range("I105").Value = 100001784
$ws.Range("J105").Value = 2149.75
$ws.Range("K105").Value = 100001784
$ws.Range("L105").Value = 2149.75
$ws.Range("M105").Value = -100000037
$ws.Range("N105").Value = -5643.75

$ws.Range("H134").Value = 7608.4736
$ws.Range("I134").Value = 968.7857
$ws.Range("K134").Value = 2906.3571
$ws.Range("M134").Value = -371.3571000000002

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 142858110
$ws.Range("I16").Value = 142858110
$ws.Range("K16").Value = 142858110
$ws.Range("M16").Value = -142857823

$ws.Range("H31").Value = 1166.3334
$ws.Range("I31").Value = 1166.3334
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1166.3334
$ws.Range("L31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -871.3334

$ws.Range("H34").Value = 1166.3334
$ws.Range("I34").Value = 1166.3334
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1166.3334
$ws.Range("L34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -964.3334

$ws.Range("H95").Value = 4860.2
$ws.Range("J95").Value = 4860.2
$ws.Range("L95").Value = 4860.2
$ws.Range("N95").Value = -10352.2

$ws.Range("H113").Value = 142858110
$ws.Range("I113").Value = 142858110
$ws.Range("K113").Value = 142858110
$ws.Range("M113").Value = -142855940

$ws.Range("H134").Value = 25642986
$ws.Range("I134").Value = 37039204
$ws.Range("J134").Value = 1500
$ws.Range("K134").Value = 111117612
$ws.Range("L134").Value = 4500
$ws.Range("M134").Value = -111115077
$ws.Range("N134").Value = -9570

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 726
$ws.Range("J16").Value = 519.2
$ws.Range("L16").Value = 519.2
$ws.Range("N16").Value = -859.2

$ws.Range("H94").Value = 13199.2
$ws.Range("J94").Value = 13199.2
$ws.Range("L94").Value = 13199.2
$ws.Range("N94").Value = -14551.2

$ws.Range("H128").Value = 35000
$ws.Range("J128").Value = 35000
$ws.Range("L128").Value = 35000
$ws.Range("N128").Value = -44960

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 10038.5
$ws.Range("J51").Value = 10038.5
$ws.Range("L51").Value = 10038.5
$ws.Range("N51").Value = -11058.5

$ws.Range("H52").Value = 10794.5
$ws.Range("J52").Value = 10794.5
$ws.Range("L52").Value = 10794.5
$ws.Range("N52").Value = -11246.5

$ws.Range("H62").Value = 125010824
$ws.Range("I62").Value = 166677330
$ws.Range("J62").Value = 11303
$ws.Range("K62").Value = 166677330
$ws.Range("L62").Value = 11303
$ws.Range("M62").Value = -166676706
$ws.Range("N62").Value = -12551

$ws.Range("H65").Value = 125010824
$ws.Range("I65").Value = 166677330
$ws.Range("J65").Value = 11303
$ws.Range("K65").Value = 833386650
$ws.Range("L65").Value = 56515
$ws.Range("M65").Value = -833383530
$ws.Range("N65").Value = -62755

$ws.Range("H113").Value = 816
$ws.Range("I113").Value = 549
$ws.Range("J113").Value = 1350
$ws.Range("K113").Value = 1647
$ws.Range("L113").Value = 4050
$ws.Range("M113").Value = 523
$ws.Range("N113").Value = -8390

$ws.Range("H123").Value = 53071.5
$ws.Range("J123").Value = 53071.5
$ws.Range("L123").Value = 53071.5
$ws.Range("N123").Value = -62871.5

$ws.Range("H126").Value = 90909990
$ws.Range("I126").Value = 90909990
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 272729970
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -272727500
